# KiCadHacks.xlsx - "more cut line moves"
# Update the MoveSegments sheet: move the cut-line text for two gr_line
# entries (rows 12/13), and refresh the dependent MIN/MAX calc rows
# (14 existing, 15 newly inserted), and renumber the trailing LEFT/MID/RIGHT
# rebuild rows that follow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MoveSegments")
$ws.Activate()

# ------------------------------------------------------------------
# 1. Offset used to grow/shrink the cut line (Z7) flips sign
# ------------------------------------------------------------------
$ws.Range("Z7").Value = 0.45

# ------------------------------------------------------------------
# 2. New gr_line source text for rows 12 & 13 (G column)
# ------------------------------------------------------------------
$ws.Range("G12").Value = "  (gr_line (start 144.78 27.559) (end 142.113 27.559) (angle 90) (layer Edge.Cuts) (width 0.127))"
$ws.Range("G13").Value = "  (gr_line (start 153.167 36.73) (end 144.533 28.097) (angle 90) (layer Edge.Cuts) (width 0.127))"

# ------------------------------------------------------------------
# 3. Row 14 formulas change (min corner) and gain a rebuilt G14 gr_line
# ------------------------------------------------------------------
$ws.Range("B14").Formula = "=MIN(B13,D13)"
$ws.Range("C14").Formula = "=MIN(C13,E13)"
$ws.Range("D14").Formula = "=MIN(B12,D12)"
$ws.Range("E14").Formula = "=C14"

$ws.Range("G16").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Formula = '=LEFT(G13,K13)&TEXT(B14,"#0.0####")&" "&TEXT(C14,"#0.0####")&MID(G13,M13-1,N13-M13+2)&TEXT(D14,"#0.0####")&" "&TEXT(E14,"#0.0####")&RIGHT(G13,LEN(G13)-P13+2)'

# ------------------------------------------------------------------
# 4. Insert a brand-new row 15 (max corner) with its own gr_line rebuild
# ------------------------------------------------------------------
$ws.Rows("15").Insert()

$ws.Range("B15").Formula = "=MAX(B13,D13)"
$ws.Range("C15").Formula = "=MIN(C13,E13)"
$ws.Range("C15").NumberFormat = "0.000"
$ws.Range("D15").Formula = "=ROUND(MAX(B12,D12)-`$Z`$7,3)"
$ws.Range("E15").Formula = "=ROUND(MAX(C12,E12)+`$Z`$7,3)"

$ws.Range("G17").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Formula = '=LEFT(G14,K13)&TEXT(B15,"#0.0####")&" "&TEXT(C15,"#0.0####")&MID(G14,M13-2,N13-M13+2)&TEXT(D15,"#0.0####")&" "&TEXT(E15,"#0.0####")&RIGHT(G14,LEN(G14)-P13+2)'

# E13 picks up the same 0.000 number format as the new C15 cell
$ws.Range("E13").NumberFormat = "0.000"

# ------------------------------------------------------------------
# 5. Rebuild rows 17-31 (the LEFT/MID/RIGHT "reconstructed gr_line" list)
#    Row 16 (old) no longer exists as a standalone entry - its formula
#    (referencing G4) now lives at row 17 after the row-15 insert; every
#    following row is re-pointed at the next source row, row 26 is gone
#    (the row that used to reference the then-empty G14 is replaced by
#    a broken #REF! formula), and the trailing shared-formula block now
#    covers G29:G31.
# ------------------------------------------------------------------
$ws.Range("G17").Formula = '=LEFT(G4,K4)&TEXT(B4,"#0.0000")&" "&TEXT(C4,"#0.0000")&MID(G4,M4-1,N4-M4+2)&TEXT(D4,"#0.0000")&" "&TEXT(E4,"#0.0000")&RIGHT(G4,LEN(G4)-P4+3)'
$ws.Range("G18").Formula = '=LEFT(G5,K5)&TEXT(B5,"#0.0000")&" "&TEXT(C5,"#0.0000")&MID(G5,M5-1,N5-M5+2)&TEXT(D5,"#0.0000")&" "&TEXT(E5,"#0.0000")&RIGHT(G5,LEN(G5)-P5+3)'
$ws.Range("G19").Formula = '=LEFT(G6,K6)&TEXT(B6,"#0.0000")&" "&TEXT(C6,"#0.0000")&MID(G6,M6-1,N6-M6+2)&TEXT(D6,"#0.0000")&" "&TEXT(E6,"#0.0000")&RIGHT(G6,LEN(G6)-P6+3)'
$ws.Range("G20").Formula = '=LEFT(G7,K7)&TEXT(B7,"#0.0000")&" "&TEXT(C7,"#0.0000")&MID(G7,M7-1,N7-M7+2)&TEXT(D7,"#0.0000")&" "&TEXT(E7,"#0.0000")&RIGHT(G7,LEN(G7)-P7+3)'
$ws.Range("G21").Formula = '=LEFT(G8,K8)&TEXT(B8,"#0.0000")&" "&TEXT(C8,"#0.0000")&MID(G8,M8-1,N8-M8+2)&TEXT(D8,"#0.0000")&" "&TEXT(E8,"#0.0000")&RIGHT(G8,LEN(G8)-P8+3)'
$ws.Range("G22").Formula = '=LEFT(G9,K9)&TEXT(B9,"#0.0000")&" "&TEXT(C9,"#0.0000")&MID(G9,M9-1,N9-M9+2)&TEXT(D9,"#0.0000")&" "&TEXT(E9,"#0.0000")&RIGHT(G9,LEN(G9)-P9+3)'
$ws.Range("G23").Formula = '=LEFT(G10,K10)&TEXT(B10,"#0.0000")&" "&TEXT(C10,"#0.0000")&MID(G10,M10-1,N10-M10+2)&TEXT(D10,"#0.0000")&" "&TEXT(E10,"#0.0000")&RIGHT(G10,LEN(G10)-P10+3)'
$ws.Range("G24").Formula = '=LEFT(G11,K11)&TEXT(B11,"#0.0000")&" "&TEXT(C11,"#0.0000")&MID(G11,M11-1,N11-M11+2)&TEXT(D11,"#0.0000")&" "&TEXT(E11,"#0.0000")&RIGHT(G11,LEN(G11)-P11+3)'
$ws.Range("G25").Formula = '=LEFT(G12,K12)&TEXT(B12,"#0.0000")&" "&TEXT(C12,"#0.0000")&MID(G12,M12-1,N12-M12+2)&TEXT(D12,"#0.0000")&" "&TEXT(E12,"#0.0000")&RIGHT(G12,LEN(G12)-P12+3)'

$ws.Rows("26").Delete()

$ws.Range("G27").Formula = '=LEFT(#REF!,K14)&TEXT(B14,"#0.0000")&" "&TEXT(C14,"#0.0000")&MID(#REF!,M14-1,N14-M14+2)&TEXT(D14,"#0.0000")&" "&TEXT(E14,"#0.0000")&RIGHT(#REF!,LEN(#REF!)-P14+3)'
$ws.Range("G28").Formula = '=LEFT(G15,K15)&TEXT(B15,"#0.0000")&" "&TEXT(C15,"#0.0000")&MID(G15,M15-1,N15-M15+2)&TEXT(D15,"#0.0000")&" "&TEXT(E15,"#0.0000")&RIGHT(G15,LEN(G15)-P15+3)'
$ws.Range("G29:G31").Formula = '=LEFT(G17,K17)&TEXT(B17,"#0.0000")&" "&TEXT(C17,"#0.0000")&MID(G17,M17-1,N17-M17+2)&TEXT(D17,"#0.0000")&" "&TEXT(E17,"#0.0000")&RIGHT(G17,LEN(G17)-P17+3)'

# ------------------------------------------------------------------
# 6. Cosmetic: widen column G a touch and autosize B:E (now populated
#    with the new min/max corner numbers), move the view / selection.
# ------------------------------------------------------------------
$ws.Columns("B:E").AutoFit() | Out-Null
$ws.Columns("G").ColumnWidth = 68.36328125

$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("G14").Select() | Out-Null

Write-Output "done"
